$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 20.5
